$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.923.82"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.665.64"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "215.54"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "0.533"
$ws.Range("E6").Value = "  +4.91%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "20.29"
$ws.Range("E10").Value = "  +3.44%  "
$ws.Range("D11").Value = "0.0898"
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("D12").Value = "1.900.21"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "1.656.07"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "66.32"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "26.899.40"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "233.71"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "8.02"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").Value = "146.16"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").Value = "1.459.82"
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("E34").Value = "  +3.76%  "
$ws.Range("D35").Value = "1.64"
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").Value = "0.577"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "0.903"
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "5.75"
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").Value = "0.979"
$ws.Range("E43").Value = "  +6.28%  "
$ws.Range("D44").Value = "65.95"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "1.808.55"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("E50").Value = "  +4.37%  "
$ws.Range("E51").Value = "  +0.54%  "
